$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "general_college_subjects.arts"
# column (R), shifting it (and everything after it) three columns to the right.
$ws.Range("R1:T1").EntireColumn.Insert()

# Give the three newly inserted header cells the same formatting
# (bold font, borders, centered alignment) as the rest of the header row.
$ws.Range("U1").Copy()
$ws.Range("R1:T1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the headers for the three newly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# The data row's new cells should default to 0, like the other
# general_college_subjects numeric columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# The "Unknown" placeholders for these particular columns are now lower-cased.
$ws.Range("D2:J2").Value = "unknown"
